# Applies the 2024-09-27 cryptos-list refresh (prices, 1h volume %, and
# the Bittensor/Binance-PegBSC-USD row swap at rows 31-32).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: write a value while forcing text storage, so numeric-looking
# strings (e.g. "1.00", "606.83") are not auto-coerced into numbers by
# Excel's usual cell-entry type inference. The NumberFormat/Style dance
# leaves no lingering style on the cell once done.
function Set-TextValue($ref, $value) {
    $rng = $ws.Range($ref)
    $rng.NumberFormat = "@"
    $rng.Value = $value
    $rng.Style = "Normal"
}

Set-TextValue "D2" '65.793.97'
Set-TextValue "D3" '2.701.18'
Set-TextValue "E3" '  +1.81%  '
Set-TextValue "E4" '  +0.04%  '
Set-TextValue "D5" '606.83'
Set-TextValue "E5" '  +1.84%  '
Set-TextValue "D6" '157.94'
Set-TextValue "E6" '  +1.00%  '
Set-TextValue "E7" '  +0.07%  '
Set-TextValue "D8" '0.588'
Set-TextValue "E8" '  -0.92%  '
Set-TextValue "E9" '  +4.89%  '
Set-TextValue "D10" '6.05'
Set-TextValue "E10" '  +4.32%  '
Set-TextValue "E11" '  +0.42%  '
Set-TextValue "E12" '  +1.20%  '
Set-TextValue "D13" '30.08'
Set-TextValue "E13" '  +3.50%  '
Set-TextValue "E14" '  +8.21%  '
Set-TextValue "D15" '3.187.92'
Set-TextValue "E15" '  +1.90%  '
Set-TextValue "D16" '65.688.25'
Set-TextValue "E16" '  +0.95%  '
Set-TextValue "D17" '2.699.84'
Set-TextValue "E17" '  +4.21%  '
Set-TextValue "D18" '12.75'
Set-TextValue "E18" '  +0.95%  '
Set-TextValue "E19" '  +1.16%  '
Set-TextValue "D20" '359.81'
Set-TextValue "E20" '  +1.34%  '
Set-TextValue "D21" '7.53'
Set-TextValue "E21" '  +3.15%  '
Set-TextValue "D22" '1.00'
Set-TextValue "E22" '  -0.12%  '
Set-TextValue "D23" '70.26'
Set-TextValue "E23" '  +2.93%  '
Set-TextValue "E24" '  +3.01%  '
Set-TextValue "E25" '  +11.58%  '
Set-TextValue "E26" '  -4.35%  '
Set-TextValue "E27" '  +2.91%  '
Set-TextValue "E28" '  +3.61%  '
Set-TextValue "D29" '8.32'
Set-TextValue "E29" '  +1.78%  '
Set-TextValue "E30" '  +4.07%  '
Set-TextValue "B31" 'Bittensor'
Set-TextValue "C31" 'https://coinranking.com/coin/pgv7xSFi6+bittensor-tao'
Set-TextValue "D31" '540.62'
Set-TextValue "E31" '  +3.44%  '
Set-TextValue "B32" 'Binance-PegBSC-USD'
Set-TextValue "C32" 'https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd'
Set-TextValue "D32" '0.994'
Set-TextValue "E32" '  -0.69%  '
Set-TextValue "E33" '  -0.27%  '
Set-TextValue "D34" '6.71'
Set-TextValue "E34" '  +5.36%  '
Set-TextValue "E35" '  -3.28%  '
Set-TextValue "E36" '  +0.75%  '
Set-TextValue "D37" '20.75'
Set-TextValue "E37" '  +2.09%  '
Set-TextValue "D38" '162.43'
Set-TextValue "E38" '  -1.54%  '
Set-TextValue "D39" '2.00'
Set-TextValue "E39" '  -1.09%  '
Set-TextValue "E40" '  -0.04%  '
Set-TextValue "D41" '0.999'
Set-TextValue "E41" '  -0.03%  '
Set-TextValue "D42" '42.82'
Set-TextValue "E42" '  +1.53%  '
Set-TextValue "D43" '168.11'
Set-TextValue "E43" '  +1.47%  '
Set-TextValue "E44" '  +1.79%  '
Set-TextValue "D45" '0.0618'
Set-TextValue "E45" '  -0.39%  '
Set-TextValue "D46" '23.60'
Set-TextValue "E46" '  +2.54%  '
Set-TextValue "E47" '  +2.30%  '
Set-TextValue "E48" '  +4.36%  '
Set-TextValue "E49" '  +1.57%  '
Set-TextValue "D50" '21.10'
Set-TextValue "E50" '  +8.03%  '
Set-TextValue "D51" '0.0985'
Set-TextValue "E51" '  -0.26%  '
